# Homework doc restructuring:
#  - Replace the "Part 1" Heading2 paragraph with a new "Overview" block
#    (section summary with time estimates) followed by a re-styled,
#    non-Heading "Section 1: Part 1: ..." paragraph.
#  - Re-style the remaining "Part 2".."Part 5" Heading2 paragraphs the
#    same way (bold, 12pt run, no paragraph style), prefixed with their
#    section number.
#  - Remove the "[Your answer here]" placeholder paragraphs.
#  - Bump the Normal style's default run size from 11pt to 12pt.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Part 1 heading -> Overview block + restyled "Section 1: ..." header
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Part 1: Sentence Type Identification") | Out-Null
$rng.Expand(4)
$xml  = '<w:p><w:pPr><w:spacing w:before="120" w:after="120"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Overview</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:ind w:left="360"/><w:spacing w:before="0" w:after="40"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Section 1: </w:t></w:r><w:r><w:t>Part 1: Sentence Type Identification (~5 min)</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:ind w:left="360"/><w:spacing w:before="0" w:after="40"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Section 2: </w:t></w:r><w:r><w:t>Part 2: Sentence Completion (~5 min)</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:ind w:left="360"/><w:spacing w:before="0" w:after="40"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Section 3: </w:t></w:r><w:r><w:t>Part 3: Sentence Writing (~5 min)</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:ind w:left="360"/><w:spacing w:before="0" w:after="40"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Section 4: </w:t></w:r><w:r><w:t>Part 4: Error Correction (~5 min)</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:ind w:left="360"/><w:spacing w:before="0" w:after="40"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Section 5: </w:t></w:r><w:r><w:t>Part 5: Analysis and Reflection (~5 min)</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:ind w:left="360"/><w:spacing w:before="120" w:after="240"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Total estimated time: </w:t></w:r><w:r><w:t>~25 minutes</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:spacing w:before="0" w:after="240"/></w:pPr><w:r><w:t>──────────────────────────────────────────────────</w:t></w:r></w:p>'
$xml += '<w:p><w:pPr><w:spacing w:before="240" w:after="120"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Section 1: Part 1: Sentence Type Identification</w:t></w:r></w:p>'
$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 2) Part 2..Part 5 headings -> restyled "Section N: ..." headers
# ---------------------------------------------------------------------
$parts = @(
  @{ Num = 2; Text = "Part 2: Sentence Completion" },
  @{ Num = 3; Text = "Part 3: Sentence Writing" },
  @{ Num = 4; Text = "Part 4: Error Correction" },
  @{ Num = 5; Text = "Part 5: Analysis and Reflection" }
)

foreach ($part in $parts) {
  $rng = $d.Content
  $rng.Find.Execute($part.Text) | Out-Null
  $rng.Expand(4)
  $newText = "Section " + $part.Num + ": " + $part.Text
  $xml = '<w:p><w:pPr><w:spacing w:before="240" w:after="120"/></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>' + $newText + '</w:t></w:r></w:p>'
  $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 3) Remove the "[Your answer here]" placeholder paragraphs
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
while ($rng.Find.Execute("[Your answer here]")) {
  $rng.Expand(4)
  $rng.Delete()
  $rng.Collapse(0)
}

# ---------------------------------------------------------------------
# 4) Normal style: bump default font size 11pt -> 12pt
# ---------------------------------------------------------------------
$d.Styles("Normal").Font.Size = 12
